{"js": "// Correct the TNS request mailbox alias: \"acstnrequest@microsoft.com\" ->\n// \"acstns@microsoft.com\". The source edit ended up with the new address\n// split across three runs (\"acstn\" / \"s\" / \"@microsoft.com\") inside the\n// existing hyperlink, so we reproduce that same run layout here.\n\nconst body = context.document.body;\n\nconst found = body.search(\"acstnrequest@microsoft.com\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  const emailRange = found.items[0];\n\n  // Replace the \"request\" portion of the address with \"s\" so the text\n  // becomes \"acstns@microsoft.com\" (formatting/hyperlink is inherited from\n  // the surrounding run).\n  const requestPart = emailRange.search(\"request\", { matchCase: true });\n  requestPart.load(\"items\");\n  await context.sync();\n  requestPart.items[0].insertText(\"s\", \"Replace\");\n  await context.sync();\n\n  // Re-locate the corrected address, then nudge (and immediately revert) a\n  // character property at each internal boundary so the text gets written\n  // back out as three discrete runs: \"acstn\" | \"s\" | \"@microsoft.com\".\n  const newAddress = body.search(\"acstns@microsoft.com\", { matchCase: true });\n  newAddress.load(\"items\");\n  await context.sync();\n  const addressRange = newAddress.items[0];\n\n  const firstSix = addressRange.search(\"acstns\", { matchCase: true });\n  firstSix.load(\"items\");\n  await context.sync();\n  const firstSixRange = firstSix.items[0];\n  firstSixRange.font.bold = true;\n  await context.sync();\n  firstSixRange.font.bold = false;\n  await context.sync();\n\n  const firstFive = addressRange.search(\"acstn\", { matchCase: true });\n  firstFive.load(\"items\");\n  await context.sync();\n  const firstFiveRange = firstFive.items[0];\n  firstFiveRange.font.bold = true;\n  await context.sync();\n  firstFiveRange.font.bold = false;\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the hyperlink run text \"acstnrequest@microsoft.com\" and change it to\n# \"acstns@microsoft.com\" (TNS alias correction), keeping it split across three\n# runs (\"acstn\" / \"s\" / \"@microsoft.com\") the way the authored edit produced.\n$r = $d.Content\n$find = $r.Find\n$find.ClearFormatting()\n$find.Text = \"acstnrequest@microsoft.com\"\n$found = $find.Execute()\n\nif ($found) {\n    $start = $r.Start\n    $end = $r.End\n\n    # \"acstnrequest@microsoft.com\"\n    #  01234567890123456789012345 6\n    # \"acstn\"(0-5) + \"request\"(5-12) + \"@microsoft.com\"(12-27)\n    # Replace the \"request\" segment with \"s\" -> \"acstn\" + \"s\" + \"@microsoft.com\"\n    $middle = $d.Range($start + 5, $start + 12)\n    $middle.Text = \"s\"\n\n    # Force the run to break into three separate runs (acstn | s | @microsoft.com)\n    # by toggling a character format on/off across each internal boundary - this\n    # mirrors how the source document ended up with three discrete <w:r> runs.\n    $boundary1 = $d.Range($start, $start + 6)\n    $boundary1.Font.Bold = 1\n    $boundary1.Font.Bold = 0\n\n    $boundary2 = $d.Range($start, $start + 5)\n    $boundary2.Font.Bold = 1\n    $boundary2.Font.Bold = 0\n}\n"}
